$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Build the two cell styles used throughout the sheet by formatting a single
# "template" cell and then propagating that exact style (format-only paste)
# to the rest of the range. This keeps the resulting styles.xml minimal
# (one shared style per distinct look) instead of one style per cell.
# ---------------------------------------------------------------------------

# Style 1: general left-aligned cell (headers, status/price/notes/counts)
$ws.Range("A1").Value = "Date"
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").Copy()
$ws.Range("B1:G9").PasteSpecial(-4122)

# Style 2: left-aligned date cell (column A data rows)
$ws.Range("A2").Value = 45650
$ws.Range("A2").HorizontalAlignment = -4131
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Price"
$ws.Range("D1").Value = "Notes"
$ws.Range("E1").Value = "MaxBookings"
$ws.Range("F1").Value = "Booked"
$ws.Range("G1").Value = "Available"

# ---------------------------------------------------------------------------
# Data rows: Date, Status, Price, Notes, MaxBookings, Booked, Available
# ---------------------------------------------------------------------------
$rows = @(
    @(2, 45650, "Limited", 18500, "Christmas Eve",  10, 7, 3),
    @(3, 45651, "Closed",  $null, "Christmas Day",   0, 0, 0),
    @(4, 45657, "Limited", 22000, "New Year's Eve", 10, 8, 2),
    @(5, 45658, "Closed",  $null, "New Year's Day",  0, 0, 0),
    @(6, 45659, "Limited", 18500, "Peak Season",    10, 6, 4),
    @(7, 45660, "Limited", 18500, "Peak Season",    10, 5, 5),
    @(8, 45698, "Closed",  $null, "Maintenance",     0, 0, 0),
    @(9, 45699, "Closed",  $null, "Maintenance",     0, 0, 0)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    if ($null -ne $row[3]) {
        $ws.Cells.Item($r, 3).Value = $row[3]
    }
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 9.7109375
$ws.Columns("B:G").ColumnWidth = 8.88671875

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$null = $ws.Range("D14").Select()

# ---------------------------------------------------------------------------
# Footer (sensitivity label stamp)
# ---------------------------------------------------------------------------
$ws.PageSetup.CenterFooter = "`r&1#&`"Calibri`"&8&K000000 Internal"
